# Coastal Surface Piercing Profilers - update Coastal CSPP cal sheet
# Correct instrument reference designators from GP05MOAS-GL003 to GP05MOAS-GL361
# and update the active sheet/selection to Asset_Cal_Info.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet: correct Ref Des and Deployment Number ---
$ws1.Range("A2").Value = "GP05MOAS-GL361"
$ws1.Range("C2").Value = 1

# --- Asset_Cal_Info sheet: correct Ref Des and Deployment Number ---
$ws2.Range("A3").Value = "GP05MOAS-GL361-00-ENG000000"
$ws2.Range("C3").Value = 1

$ws2.Range("A4").Value = "GP05MOAS-GL361-01-FLORDM000"
$ws2.Range("C4").Value = 1

$ws2.Range("A5").Value = "GP05MOAS-GL361-01-FLORDM000"
$ws2.Range("C5").Value = 1

$ws2.Range("A6").Value = "GP05MOAS-GL361-01-FLORDM000"
$ws2.Range("C6").Value = 1

$ws2.Range("A7").Value = "GP05MOAS-GL361-01-FLORDM000"
$ws2.Range("C7").Value = 1

$ws2.Range("A8").Value = "GP05MOAS-GL361-02-DOSTAM000"
$ws2.Range("C8").Value = 1

$ws2.Range("A9").Value = "GP05MOAS-GL361-04-CTDGVM000"
$ws2.Range("C9").Value = 1

# --- Update selections on each sheet ---
$ws1.Activate()
$ws1.Range("D11:D12").Select()

# Asset_Cal_Info becomes the active/selected tab
$ws2.Activate()
$ws2.Range("C27").Select()
